$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C4").Value = 50000
$wb.Application.Calculate()

$ws.Columns.Item(4).AutoFit() | Out-Null
$ws.Columns.Item(15).AutoFit() | Out-Null
